$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1837349397590362
$ws.Range("C2").Value = 0.5903614457831325
$ws.Range("J2").Value = 0.006024096385542169
$ws.Range("P2").Value = 0.1325301204819277
$ws.Range("S2").Value = 0.08734939759036145
$ws.Range("C3").Value = 0.025
$ws.Range("J3").Value = 0.03
$ws.Range("P3").Value = 0.73
$ws.Range("S3").Value = 0.215
$ws.Range("J4").Value = 0.02631578947368421
$ws.Range("P4").Value = 0.631578947368421
$ws.Range("S4").Value = 0.3421052631578947
$ws.Range("B6").Value = 0.06130268199233716
$ws.Range("D6").Value = 0.01149425287356322
$ws.Range("F6").Value = 0.04980842911877394
$ws.Range("J6").Value = 0.3984674329501915
$ws.Range("O6").Value = 0.02298850574712644
$ws.Range("Q6").Value = 0.1226053639846743
$ws.Range("R6").Value = 0.03065134099616858
$ws.Range("S6").Value = 0.3026819923371648
$ws.Range("B7").Value = 0.1158301158301158
$ws.Range("D7").Value = 0.01158301158301158
$ws.Range("F7").Value = 0.04633204633204633
$ws.Range("J7").Value = 0.1505791505791506
$ws.Range("O7").Value = 0.03474903474903475
$ws.Range("Q7").Value = 0.1853281853281853
$ws.Range("R7").Value = 0.03861003861003861
$ws.Range("S7").Value = 0.416988416988417
$ws.Range("B8").Value = 0.09829059829059829
$ws.Range("D8").Value = 0.02564102564102564
$ws.Range("F8").Value = 0.05982905982905983
$ws.Range("J8").Value = 0.1068376068376068
$ws.Range("O8").Value = 0.01282051282051282
$ws.Range("Q8").Value = 0.1623931623931624
$ws.Range("R8").Value = 0.05128205128205128
$ws.Range("S8").Value = 0.4829059829059829
$ws.Range("B9").Value = 0.1341463414634146
$ws.Range("D9").Value = 0.006097560975609756
$ws.Range("F9").Value = 0.03658536585365853
$ws.Range("J9").Value = 0.1097560975609756
$ws.Range("O9").Value = 0.03048780487804878
$ws.Range("Q9").Value = 0.1768292682926829
$ws.Range("R9").Value = 0.07317073170731707
$ws.Range("S9").Value = 0.4329268292682927
$ws.Range("B10").Value = 0.1236777868185517
$ws.Range("D10").Value = 0.01790073230268511
$ws.Range("F10").Value = 0.06346623270951994
$ws.Range("J10").Value = 0.1106590724165989
$ws.Range("O10").Value = 0.01952807160292921
$ws.Range("Q10").Value = 0.2351505288852726
$ws.Range("R10").Value = 0.04963384865744508
$ws.Range("S10").Value = 0.3799837266069976
$ws.Range("G11").Value = 0.1401869158878505
$ws.Range("J11").Value = 0.1168224299065421
$ws.Range("K11").Value = 0.2079439252336449
$ws.Range("L11").Value = 0.5163551401869159
$ws.Range("S11").Value = 0.01869158878504673
$ws.Range("F12").Value = 0.004444444444444444
$ws.Range("G12").Value = 0.7377777777777778
$ws.Range("J12").Value = 0.1911111111111111
$ws.Range("K12").Value = 0.008888888888888889
$ws.Range("L12").Value = 0.02222222222222222
$ws.Range("S12").Value = 0.03555555555555556
$ws.Range("F13").Value = 0.01886792452830189
$ws.Range("G13").Value = 0.7358490566037735
$ws.Range("J13").Value = 0.1509433962264151
$ws.Range("S13").Value = 0.09433962264150944
$ws.Range("F15").Value = 0.03896103896103896
$ws.Range("H15").Value = 0.1818181818181818
$ws.Range("I15").Value = 0.0735930735930736
$ws.Range("J15").Value = 0.329004329004329
$ws.Range("K15").Value = 0.08658008658008658
$ws.Range("M15").Value = 0.01731601731601732
$ws.Range("N15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.08225108225108226
$ws.Range("S15").Value = 0.1861471861471861
$ws.Range("F16").Value = 0.0673076923076923
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.05288461538461538
$ws.Range("J16").Value = 0.4038461538461539
$ws.Range("K16").Value = 0.125
$ws.Range("M16").Value = 0.01923076923076923
$ws.Range("O16").Value = 0.07692307692307693
$ws.Range("S16").Value = 0.1298076923076923
$ws.Range("F17").Value = 0.05106382978723404
$ws.Range("H17").Value = 0.1978723404255319
$ws.Range("I17").Value = 0.06808510638297872
$ws.Range("J17").Value = 0.3319148936170213
$ws.Range("K17").Value = 0.1617021276595745
$ws.Range("M17").Value = 0.0276595744680851
$ws.Range("O17").Value = 0.0574468085106383
$ws.Range("S17").Value = 0.1042553191489362
$ws.Range("F18").Value = 0.05084745762711865
$ws.Range("H18").Value = 0.1949152542372881
$ws.Range("I18").Value = 0.05084745762711865
$ws.Range("J18").Value = 0.3898305084745763
$ws.Range("K18").Value = 0.1355932203389831
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.06779661016949153
$ws.Range("S18").Value = 0.09322033898305085
$ws.Range("F19").Value = 0.03933434190620273
$ws.Range("H19").Value = 0.2186081694402421
$ws.Range("I19").Value = 0.07110438729198185
$ws.Range("J19").Value = 0.3305597579425114
$ws.Range("K19").Value = 0.1391830559757943
$ws.Range("M19").Value = 0.02344931921331316
$ws.Range("O19").Value = 0.06127080181543117
$ws.Range("S19").Value = 0.1164901664145234
